$d = $word.ActiveDocument

# Locate the "Known bugs" bullet paragraph that currently reads:
# "When the user receives a notification, the user cannot click in the
#  menu as long as the notification is present"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "When the user receives a notification*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$full = $target.Range   # paragraph range, includes the trailing paragraph mark

# Grab the original opening <w:p ...> tag (with its paraId/rsid attributes)
# so the first paragraph keeps its identity after the rewrite.
$ooxml = $full.WordOpenXML
$startIdx = $ooxml.IndexOf("<w:p ")
$endIdx = $ooxml.IndexOf(">", $startIdx)
$pOpenTag = $ooxml.Substring($startIdx, $endIdx - $startIdx + 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
$pOpenTag +
'<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">When the user receives a notification, the user cannot click in the menu </w:t></w:r>' +
'<w:proofErr w:type="gramStart"/>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>as long as</w:t></w:r>' +
'<w:proofErr w:type="gramEnd"/>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the notification is present</w:t></w:r>' +
'<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
'</w:p>' +
'<w:p>' +
'<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The hamburger menu is too much on the right side for mobile phones</w:t></w:r>' +
'</w:p>' +
'<w:p>' +
'<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Log out does not work</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml) | Out-Null
